$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (row 2 through row 10) from serial date 45204 to 45207
for ($r = 2; $r -le 10; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45207
    }
}
